# Commit: "Fruta / hortaliza, semanal"
# A new weekly data row is inserted at row 190 (pushing the existing
# rows 190-309 down to 191-310); the new row contains a fresh
# observation for Perejil - Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 190; everything below shifts down
# by one (old row 190 becomes 191, ..., old row 309 becomes 310).
$ws.Rows(190).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A190").Value = 9
$ws.Range("B190").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C190").Value = "Metropolitana"
$ws.Range("D190").Value = 44606
$ws.Range("E190").Value = 13
$ws.Range("F190").Value = 100112044
$ws.Range("G190").Value = "Perejil"
$ws.Range("H190").Value = "Sin especificar"
$ws.Range("I190").Value = "Primera"
$ws.Range("J190").Value = 43
$ws.Range("K190").Value = 16000
$ws.Range("L190").Value = 18000
$ws.Range("M190").Value = 17023
$ws.Range("N190").Value = "$/docena de atados"
$ws.Range("O190").Value = "Región Metropolitana"
$ws.Range("P190").Value = 5674
$ws.Range("Q190").Value = 3
$ws.Range("R190").Value = "Hortaliza"
